$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leading index column (A) - this shifts B:F left to A:E,
# turning the former B..F header row into A..E and dropping the old
# index values in column A (1 / 8) entirely.
$ws.Columns("A").Delete()

# The header that used to read "MODEL_CONDITION" (now in column D after
# the shift) is renamed to "MODELCONDITION".
$ws.Range("D1").Value = "MODELCONDITION"
